$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 76, shifting existing rows 76-159 down to 77-160.
$ws.Rows.Item(76).Insert()

# Populate the newly inserted row 76 with the new data record.
$ws.Range("A76").Value = 4
$ws.Range("B76").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C76").Value = "Los Lagos"
$ws.Range("D76").Value = 44539
$ws.Range("E76").Value = 10
$ws.Range("F76").Value = "Fruta"
$ws.Range("G76").Value = 100109
$ws.Range("H76").Value = "Uva"
$ws.Range("I76").Value = 100109001
$ws.Range("J76").Value = "Uva"
$ws.Range("K76").Value = "Superior Seedless"
$ws.Range("L76").Value = "Primera"
$ws.Range("M76").Value = 200
$ws.Range("N76").Value = 26000
$ws.Range("O76").Value = 27000
$ws.Range("P76").Value = 26500
$ws.Range("Q76").Value = "`$/bandeja 8 kilos"
$ws.Range("R76").Value = "Provincia de Limarí"
$ws.Range("S76").Value = 3312
$ws.Range("T76").Value = 8
